# Adds a new "Feb 21" timesheet section (rows 183-190) to the existing
# Timesheet worksheet: a blank separator row followed by seven entries
# covering 10:00-19:00 (with a lunch break and a client call), matching
# the formatting of the preceding day-sections already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Blank separator row, formatted like the other day separators (row 20) ---
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A183:C183").PasteSpecial(-4122) | Out-Null

# --- 2. First entry of the day: taller "description" row style (like row 4) ---
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A184:C184").PasteSpecial(-4122) | Out-Null

# --- 3. Middle entries: plain entry row style (like row 3) ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A185:C185").PasteSpecial(-4122) | Out-Null
$ws.Range("A186:C186").PasteSpecial(-4122) | Out-Null
$ws.Range("A187:C187").PasteSpecial(-4122) | Out-Null
$ws.Range("A188:C188").PasteSpecial(-4122) | Out-Null
$ws.Range("A189:C189").PasteSpecial(-4122) | Out-Null

# --- 4. Last entry of the day: taller "description" row style again (like row 4) ---
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A190:C190").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 5. Fill in the values ---
$ws.Range("A184").Value = "Feb 21 10:00 to 11:00"
$ws.Range("B184").Value = "Builded Decision tree, random forest, support vector machine and k-nearest neighbour`nmodels and compared there accuracy."
$ws.Range("C184").Value = "Infimetrics"

$ws.Range("A185").Value = "Feb 21 11:00 to 12:00"
$ws.Range("B185").Value = "Builded multilayer perceptron model, modified some code. Working on model tuning"
$ws.Range("C185").Value = "Infimetrics"

$ws.Range("A186").Value = "Feb 21 12:00 to 13:00"
$ws.Range("B186").Value = "Client call"
$ws.Range("C186").Value = "Sapphire auomation"

$ws.Range("A187").Value = "Feb 21 13:00 to 14:00"
$ws.Range("B187").Value = "Practiced some examples of django based deployement of ml models"
$ws.Range("C187").Value = "Infimetrics"

$ws.Range("A188").Value = "Feb 21 14:00 to 15:00"
$ws.Range("B188").Value = "Lunch"
$ws.Range("C188").Value = "Infimetrics"

$ws.Range("A189").Value = "Feb 21 15:00 to 16:00"
$ws.Range("B189").Value = "Modified code. Added code which save model as pickled model"
$ws.Range("C189").Value = "Infimetrics"

$ws.Range("A190").Value = "Feb 21 16:00 to 19:00"
$ws.Range("B190").Value = "Worked on many django based ml model deployement examples, none were executed`nsuccessfully."
$ws.Range("C190").Value = "Infimetrics"

# --- 6. These two rows hold a long wrapped description that needs three
#        visual lines (matching the other multi-line "description" rows in
#        the sheet, e.g. row 4/182) ---
$ws.Rows.Item(184).RowHeight = 45
$ws.Rows.Item(190).RowHeight = 45

# --- 7. Update the view so the new last row is visible/selected, same as the
#        author's save (scrolled down, D190 selected) ---
$ws.Range("D190").Select()
